$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.161.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.850.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4711"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2899"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06537"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07945"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.863.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.079"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6761"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "269.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.141.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007516"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.101.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.219"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.123"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "165.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.134"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.930"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.394"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09863"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.276"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.991"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04677"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6971"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01868"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.607"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.309"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.925"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8353"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4136"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "944.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.119"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.963"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05647"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.30%  "
